$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.746.81"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "1.602.80"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("E10").Value = "  +0.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0847"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "1.827.71"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").Value = "1.609.14"
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.95%  "
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("E22").Value = "  -3.64%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("E26").Value = "  -0.66%  "
$ws.Range("E27").Value = "  -0.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("E32").Value = "  +0.41%  "
$ws.Range("D33").Value = "1.290.72"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("E36").Value = "  -2.34%  "
$ws.Range("E37").Value = "  +12.18%  "
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.832"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.782"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "63.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("D44").Value = "1.739.82"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("E45").Value = "  -0.36%  "
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0104"
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.102"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0516"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.68%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.33%  "
